$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 19: Aula 30 - Template baseado em layout
$ws.Cells.Item(19, 2).Value = 30
$ws.Cells.Item(19, 3).Value = "7. Thymeleaf para as Views"
$ws.Cells.Item(19, 4).Value = "30. Template baseado em layout"

# Observation text (E19) has two differently-formatted runs: plain text,
# then a bold/red/Calibri-11 run with the xmlns snippet.
$run1 = "2:50 - para criar templeate baseado em layout, o template modelo HTML do thymeleaf deve conter os namespaces necessários a seguir;`n"
$run2 = "`nxmlns:th=`"http://www.thymeleaf.org`"`nxmlns:layout=`"http://www.ultraq.net.nz/thymeleaf/layout`""
$fullText = $run1 + $run2

$ws.Cells.Item(19, 5).Value = $fullText
$ws.Cells.Item(19, 5).WrapText = $true

$rng = $ws.Cells.Item(19, 5)
$startPos = $run1.Length + 1
$runLen = $run2.Length
$chars = $rng.Characters($startPos, $runLen)
$chars.Font.Bold = $true
$chars.Font.Color = 255
$chars.Font.Name = "Calibri"
$chars.Font.Size = 11

# Row height matches the other wrapped-text rows (75pt, same as row 5).
$ws.Range("B19:E19").RowHeight = 75

# Matches the author's final selection recorded in the sheet view.
$ws.Range("E23:E26").Select()
